$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 175.5038511461806
$ws.Range("C2").Value = 246.099340738569
$ws.Range("D2").Value = 283.4982903894888
$ws.Range("E2").Value = 306.3472503517436

$ws.Range("B3").Value = 188.077694752631
$ws.Range("C3").Value = 267.5086431107969
$ws.Range("D3").Value = 312.677283613454
$ws.Range("E3").Value = 346.7295360882136

$ws.Range("B4").Value = 112.4248170957789
$ws.Range("C4").Value = 166.9294644886803
$ws.Range("D4").Value = 218.1011528086006
$ws.Range("E4").Value = 257.7098228245226

$ws.Range("B5").Value = 58.71294989186963
$ws.Range("C5").Value = 69.77250912507837
$ws.Range("D5").Value = 68.37877842506221
$ws.Range("E5").Value = 64.47725878194382

$ws.Range("B6").Value = 41.8105744583499
$ws.Range("C6").Value = 44.24276505268583
$ws.Range("D6").Value = 40.57225876513246
$ws.Range("E6").Value = 33.93938830210805

$ws.Range("B7").Value = 4.929144550771079
$ws.Range("C7").Value = 5.763203358252913
$ws.Range("D7").Value = 6.023653289621835
$ws.Range("E7").Value = 5.559129768194828

$ws.Range("B8").Value = 159.2734554290184
$ws.Range("C8").Value = 232.3592039582677
$ws.Range("D8").Value = 250.02917500901
$ws.Range("E8").Value = 253.942977826177

$ws.Range("B9").Value = 116.869797503889
$ws.Range("C9").Value = 135.7545564028199
$ws.Range("D9").Value = 116.2066623531116
$ws.Range("E9").Value = 98.01103415991363

$ws.Range("B10").Value = 37.13227483634517
$ws.Range("C10").Value = 60.66173511600021
$ws.Range("D10").Value = 88.00652696885979
$ws.Range("E10").Value = 109.6781906732005

$ws.Range("B11").Value = 5.993161529336554
$ws.Range("C11").Value = 9.173977692284392
$ws.Range("D11").Value = 13.53124367337835
$ws.Range("E11").Value = 18.60696974524193

$ws.Range("B12").Value = 3.590279878239165
$ws.Range("C12").Value = 7.006436624380287
$ws.Range("D12").Value = 15.66120373963864
$ws.Range("E12").Value = 23.97667024404713

$ws.Range("B13").Value = 15.92651745250431
$ws.Range("C13").Value = 26.22543069112396
$ws.Range("D13").Value = 41.40487310090964
$ws.Range("E13").Value = 54.65597690983485

